$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray county-name header row (K28:AD28) and the three
# panel-data rows right under it (29:31) -- the "still fail" attempt
# whose numbers are being pulled back out.
$ws.Range("K28:AD28").ClearContents()
$ws.Range("J29:AD31").ClearContents()

# Also blank out the duplicated panel block (AO37:BH39); keep the cell
# formatting/borders in place, only the values go.
$ws.Range("AO37:BH39").ClearContents()

# Scroll/selection housekeeping: drop the pinned topLeftCell and leave
# the selection sitting over the now-empty block.
$ws.Activate()
$ws.Range("AO36:BI40").Select()
